$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Preserve B2's current (Hyperlink) cell format before touching hyperlinks ---
$ws.Range("B2").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# --- Remove hyperlinks from B3 and B4 (engine only supports deleting all at once) ---
$ws.Hyperlinks.Delete()

# --- Re-add the hyperlink that must remain, on B2 ---
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:webuser@unilogcorp.com")

# --- Restore B2's original formatting (Hyperlinks.Add() can reset it) ---
$ws.Range("H1").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("H1").Clear()

# --- B4: switch away from the Hyperlink style to a plain quote-prefixed / numberformat style ---
$ws.Range("A2").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("B4").NumberFormat = "0"

# --- Update the changed cell values ---
$ws.Range("B4").Value = "'generaluser2@unilogcorp.com"
$ws.Range("D4").Value = "General User Test User"

# --- Update the stored selection ---
$ws.Range("G9").Select()
